$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, pushing existing rows 138:171 down to 139:172
$ws.Rows("138:138").Insert()

# Populate the new row 138 with the new record's data
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C138").Value = 'Los Lagos'
$ws.Range("D138").Value = 44511
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 100112003
$ws.Range("G138").Value = 'Ajo'
$ws.Range("H138").Value = 'Chino'
$ws.Range("I138").Value = 'Primera'
$ws.Range("J138").Value = 80
$ws.Range("K138").Value = 21000
$ws.Range("L138").Value = 22000
$ws.Range("M138").Value = 21500
$ws.Range("N138").Value = '$/caja 10 kilos'
$ws.Range("O138").Value = 'China'
$ws.Range("P138").Value = 2150
$ws.Range("Q138").Value = 10
$ws.Range("R138").Value = 'Hortaliza'
